# Apprentices.xlsx refactor ("fixed view add new icon and a lot of refactoring"):
#  - fill in the header that the view was missing: A1 ("Name")
#  - add a new apprentice entry as row 3 (Name / Email / Phone / Project name /
#    Occupation / Added date)
#  - refresh column widths so they best-fit the new (longer) content, same as
#    Format > AutoFit Column Width in the UI

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the view: the header row was missing its first label.
$ws.Range("A1").Value = "Name"

# New apprentice entry (row 3). Name/Email/Phone/Occupation are stored
# encrypted in the source data (long opaque text), same shape as the
# existing PII already in row 2.
$ws.Range("A3").Value = "NBV6qZjWZyYdi1kAKln+46n07FtdDKgPUaweBSoWUdDlqqPbKDPFiVtShe0dg60wXNfW5+w0O9sUi/gflkooQh+998aT7iMKOMuYQhTlX89A2mGRO/mH3pXljHU4BubEPP4YZqrW+9qCA8m+RtzjXpOCDeaIe99kOQQibe0H5gw="
$ws.Range("B3").Value = "bfXCkN6OI/k6lKWAlmqlxuZKUYjC2v+ppzktjsx7WJsblgWbdwmWOq8at1DIGO3kZBqLqbhgzuQQyGdNJMaiVr3r3ZrlumL+y0snAonQpS9CvdUhgbCghHQy8ktzp4pw11GsQoGTxQFI/hOvr/SsGhZiBSWteuJt/KceDVJKN/I="
$ws.Range("C3").Value = "d+eFclghqlZ51NbGIVZg2UhdyCNJj0Mo+MvGLF6ry1vUEgcPtJuuLfdsGdQ8L2XJWjlxBqWlsDwm7p8JOBoHbm+OHQfBMA582n3NF2irdMFwwsKQ6ja9D7dFBexBOJgWnJWlZlaP/cOzhdJbvVY/GGxu8q/jbqmk7ZGG3IbpO3s="
$ws.Range("D3").Value = "CLS"
$ws.Range("E3").Value = "WS5sTbLXc36YJqMXyGb+souFJzoMuKvRnkr0gZ9b+mEw0hg0QVVpEPp7/usNy3OWD29YcZri7SpvuezsDXwAm/sQ4lTdXywp777wNGNBwSZtLYlvT+Q6IxnP+vbVllijWQyedMmjS6hYGG47fHFzJEfEcHO2UF+Zge31345kM+Y="
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2017-05-23"

# Re-fit the columns to the new, wider content.
$ws.Columns.Item(1).EntireColumn.AutoFit()
$ws.Columns.Item(2).EntireColumn.AutoFit()
$ws.Columns.Item(3).EntireColumn.AutoFit()
$ws.Columns.Item(4).EntireColumn.AutoFit()
$ws.Columns.Item(5).EntireColumn.AutoFit()
$ws.Columns.Item(6).EntireColumn.AutoFit()
